$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ln_real_gdp
$ws.Range("B2").Value = 6.438261300782886
$ws.Range("C2").Value = 6.501110433718534
$ws.Range("D2").Value = 8.145282382069622
$ws.Range("E2").Value = "Yes"

# Row 3: ln_人口密度
$ws.Range("B3").Value = 5.722068186375838
$ws.Range("C3").Value = 5.768180621157776
$ws.Range("D3").Value = 5.370774989784516

# Row 4: ln_金融发展水平
$ws.Range("B4").Value = 0.6021899870632997
$ws.Range("C4").Value = 0.5812001176826126
$ws.Range("D4").Value = 6.398799704345692

# Row 5: 第二产业占GDP比重
$ws.Range("B5").Value = 0.497744963449168
$ws.Range("C5").Value = 0.4903009433128748
$ws.Range("D5").Value = 7.43668508647872
$ws.Range("E5").Value = "Yes"
